$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 10000
$ws.Range("I48").Value = 10000
$ws.Range("K48").Value = 30000
$ws.Range("M48").Value = -29708
$ws.Range("H56").Value = 10000
$ws.Range("I56").Value = 10000
$ws.Range("K56").Value = 30000
$ws.Range("M56").Value = -29466
$ws.Range("H64").Value = 31255160
$ws.Range("J64").Value = 5811.846
$ws.Range("L64").Value = 5811.846
$ws.Range("N64").Value = -6307.846
$ws.Range("H67").Value = 31255160
$ws.Range("J67").Value = 5811.846
$ws.Range("L67").Value = 5811.846
$ws.Range("N67").Value = -7527.846
$ws.Range("H106").Value = 5458.1665
$ws.Range("I106").Value = 5458.1665
$ws.Range("K106").Value = 5458.1665
$ws.Range("M106").Value = -4827.1665

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 2118857.8
$ws.Range("I5").Value = 5650150
$ws.Range("K5").Value = 5650150
$ws.Range("M5").Value = -5650038
$ws.Range("H32").Value = 5322345.5
$ws.Range("I32").Value = 5683983
$ws.Range("K32").Value = 5683983
$ws.Range("M32").Value = -5683696
$ws.Range("H45").Value = 3292.5715
$ws.Range("I45").Value = 3341.3333
$ws.Range("K45").Value = 3341.3333
$ws.Range("M45").Value = -2964.3333
$ws.Range("H61").Value = 2117493.5
$ws.Range("I61").Value = 2932.5
$ws.Range("J61").Value = 3527200.8
$ws.Range("K61").Value = 2932.5
$ws.Range("L61").Value = 3527200.8
$ws.Range("M61").Value = -2720.5
$ws.Range("N61").Value = -3527624.8
$ws.Range("H74").Value = 13882.892
$ws.Range("I74").Value = 1081.359
$ws.Range("K74").Value = 1081.359
$ws.Range("M74").Value = -207.3589999999999
$ws.Range("H77").Value = 13882.892
$ws.Range("I77").Value = 1081.359
$ws.Range("K77").Value = 5406.795
$ws.Range("M77").Value = -1038.795
$ws.Range("H122").Value = 769646.3
$ws.Range("I122").Value = 1012371.7
$ws.Range("K122").Value = 3037115.1
$ws.Range("M122").Value = -3034665.1
$ws.Range("H136").Value = 2117493.5
$ws.Range("I136").Value = 2932.5
$ws.Range("J136").Value = 3527200.8
$ws.Range("K136").Value = 8797.5
$ws.Range("L136").Value = 10581602.4
$ws.Range("M136").Value = -6247.5
$ws.Range("N136").Value = -10586702.4
$ws.Range("H140").Value = 76666.336
$ws.Range("J140").Value = 76666.336
$ws.Range("L140").Value = 76666.336
$ws.Range("N140").Value = -87026.336

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 2118857.8
$ws.Range("I4").Value = 5650150
$ws.Range("K4").Value = 5650150
$ws.Range("M4").Value = -5650035
$ws.Range("H99").Value = 18254.477
$ws.Range("I99").Value = 19465.475
$ws.Range("J99").Value = 6750
$ws.Range("K99").Value = 19465.475
$ws.Range("L99").Value = 6750
$ws.Range("M99").Value = -17967.475
$ws.Range("N99").Value = -9746
$ws.Range("H107").Value = 1380.1034
$ws.Range("I107").Value = 1347.25
$ws.Range("K107").Value = 1347.25
$ws.Range("M107").Value = 572.75
$ws.Range("H134").Value = 54181.793
$ws.Range("I134").Value = 80141.16
$ws.Range("K134").Value = 240423.48
$ws.Range("M134").Value = -237888.48

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H37").Value = 7161.4
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 7161.4
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 7161.4
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -7375.4
$ws.Range("H86").Value = 12893.77
$ws.Range("J86").Value = 6244.1665
$ws.Range("L86").Value = 6244.1665
$ws.Range("N86").Value = -8490.166499999999
$ws.Range("H89").Value = 12893.77
$ws.Range("J89").Value = 6244.1665
$ws.Range("L89").Value = 31220.8325
$ws.Range("N89").Value = -42452.8325
$ws.Range("H97").Value = 17000
$ws.Range("J97").Value = 17000
$ws.Range("L97").Value = 17000
$ws.Range("N97").Value = -18982
$ws.Range("H132").Value = 32171616
$ws.Range("I132").Value = 2182.7036
$ws.Range("K132").Value = 6548.110799999999
$ws.Range("M132").Value = -4018.110799999999
$ws.Range("H134").Value = 26318544
$ws.Range("I134").Value = 1334.3478
$ws.Range("J134").Value = 66671600
$ws.Range("K134").Value = 4003.0434
$ws.Range("L134").Value = 200014800
$ws.Range("M134").Value = -1468.0434
$ws.Range("N134").Value = -200019870

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 2532.5
$ws.Range("J26").Value = 60
$ws.Range("L26").Value = 180
$ws.Range("N26").Value = -756
$ws.Range("H81").Value = 17001000
$ws.Range("J81").Value = 17001000
$ws.Range("L81").Value = 51003000
$ws.Range("N81").Value = -51005246
$ws.Range("H84").Value = 17001000
$ws.Range("J84").Value = 17001000
$ws.Range("L84").Value = 153009000
$ws.Range("N84").Value = -153020232

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1113371.2
$ws.Range("I40").Value = 2504.9033
$ws.Range("K40").Value = 2504.9033
$ws.Range("M40").Value = -2368.9033
$ws.Range("H68").Value = 50000
$ws.Range("I68").Value = 50000
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 50000
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -49251
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 50000
$ws.Range("I71").Value = 50000
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 250000
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -246256
$ws.Range("N71").ClearContents()
$ws.Range("H100").Value = 3536
$ws.Range("I100").Value = 2950.6
$ws.Range("J100").Value = 4999.5
$ws.Range("K100").Value = 2950.6
$ws.Range("L100").Value = 4999.5
$ws.Range("M100").Value = -2409.6
$ws.Range("N100").Value = -6081.5
$ws.Range("H132").Value = 1493100.6
$ws.Range("I132").Value = 2931.9524
$ws.Range("K132").Value = 8795.8572
$ws.Range("M132").Value = -6265.8572
$ws.Range("H136").Value = 1506582
$ws.Range("J136").Value = 2719350.8
$ws.Range("L136").Value = 8158052.399999999
$ws.Range("N136").Value = -8163152.399999999

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 6000
$ws.Range("I39").Value = 6000
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 6000
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -5587
$ws.Range("N39").ClearContents()
$ws.Range("H43").Value = 25000
$ws.Range("I43").Value = 15000
$ws.Range("K43").Value = 15000
$ws.Range("M43").Value = -14851
$ws.Range("H96").Value = 1190.931
$ws.Range("I96").Value = 990.8261
$ws.Range("K96").Value = 990.8261
$ws.Range("M96").Value = 382.1739
$ws.Range("H132").Value = 368340.4
$ws.Range("I132").Value = 3466.4211
$ws.Range("K132").Value = 10399.2633
$ws.Range("M132").Value = -7869.263300000001
$ws.Range("H136").Value = 247057.67
$ws.Range("J136").Value = 660906.75
$ws.Range("L136").Value = 1982720.25
$ws.Range("N136").Value = -1987820.25
